$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(1)
$sh.TextFrame.TextRange.Font.Bold = $false
